$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(131, 8).Value = 1777.7812
$ws.Cells.Item(131, 9).Value = 690.8261
$ws.Cells.Item(131, 10).Value = 4555.5557
$ws.Cells.Item(131, 11).Value = 2072.4783
$ws.Cells.Item(131, 12).Value = 13666.6671
$ws.Cells.Item(131, 13).Value = 2967.5217
$ws.Cells.Item(131, 14).Value = -23746.6671
$ws.Cells.Item(135, 8).Value = 12928337
$ws.Cells.Item(135, 9).Value = 432.39285
$ws.Cells.Item(135, 10).Value = 40773052
$ws.Cells.Item(135, 11).Value = 3891.53565
$ws.Cells.Item(135, 12).Value = 366957468
$ws.Cells.Item(135, 13).Value = -1356.53565
$ws.Cells.Item(135, 14).Value = -366962538
$ws.Cells.Item(137, 8).Value = 20001186
$ws.Cells.Item(137, 9).Value = 1112.3235
$ws.Cells.Item(137, 10).Value = 62501344
$ws.Cells.Item(137, 11).Value = 3336.9705
$ws.Cells.Item(137, 12).Value = 187504032
$ws.Cells.Item(137, 13).Value = -786.9704999999999
$ws.Cells.Item(137, 14).Value = -187509132
$ws.Cells.Item(138, 8).Value = 2792.3809
$ws.Cells.Item(138, 9).Value = 2406.0715
$ws.Cells.Item(138, 10).Value = 3101.4285
$ws.Cells.Item(138, 11).Value = 7218.2145
$ws.Cells.Item(138, 12).Value = 9304.2855
$ws.Cells.Item(138, 13).Value = -2078.2145
$ws.Cells.Item(138, 14).Value = -19584.2855
$ws.Cells.Item(141, 8).Value = 1085.3962
$ws.Cells.Item(141, 9).Value = 454.65
$ws.Cells.Item(141, 10).Value = 3026.1538
$ws.Cells.Item(141, 11).Value = 1363.95
$ws.Cells.Item(141, 12).Value = 9078.4614
$ws.Cells.Item(141, 13).Value = 3816.05
$ws.Cells.Item(141, 14).Value = -19438.4614
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1473.99
$ws.Cells.Item(32, 9).Value = 969.4524
$ws.Cells.Item(32, 10).Value = 4122.8125
$ws.Cells.Item(32, 11).Value = 969.4524
$ws.Cells.Item(32, 12).Value = 4122.8125
$ws.Cells.Item(32, 13).Value = -682.4524
$ws.Cells.Item(32, 14).Value = -4696.8125
$ws.Cells.Item(61, 8).Value = 2268706.2
$ws.Cells.Item(61, 9).Value = 2646615.5
$ws.Cells.Item(61, 10).Value = 1250
$ws.Cells.Item(61, 11).Value = 2646615.5
$ws.Cells.Item(61, 12).Value = 1250
$ws.Cells.Item(61, 13).Value = -2646403.5
$ws.Cells.Item(61, 14).Value = -1674
$ws.Cells.Item(74, 8).Value = 10642154
$ws.Cells.Item(74, 9).Value = 13514210
$ws.Cells.Item(74, 11).Value = 13514210
$ws.Cells.Item(74, 13).Value = -13513336
$ws.Cells.Item(77, 8).Value = 10642154
$ws.Cells.Item(77, 9).Value = 13514210
$ws.Cells.Item(77, 11).Value = 67571050
$ws.Cells.Item(77, 13).Value = -67566682
$ws.Cells.Item(132, 8).Value = 6252629
$ws.Cells.Item(132, 9).Value = 7275983.5
$ws.Cells.Item(132, 10).Value = 112500.445
$ws.Cells.Item(132, 11).Value = 21827950.5
$ws.Cells.Item(132, 12).Value = 337501.335
$ws.Cells.Item(132, 13).Value = -21825420.5
$ws.Cells.Item(132, 14).Value = -342561.335
$ws.Cells.Item(136, 8).Value = 2268706.2
$ws.Cells.Item(136, 9).Value = 2646615.5
$ws.Cells.Item(136, 10).Value = 1250
$ws.Cells.Item(136, 11).Value = 7939846.5
$ws.Cells.Item(136, 12).Value = 3750
$ws.Cells.Item(136, 13).Value = -7937296.5
$ws.Cells.Item(136, 14).Value = -8850
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 282.44446
$ws.Cells.Item(80, 9).Value = 258.2
$ws.Cells.Item(80, 10).Value = 291.76923
$ws.Cells.Item(80, 11).Value = 258.2
$ws.Cells.Item(80, 12).Value = 291.76923
$ws.Cells.Item(80, 13).Value = 739.8
$ws.Cells.Item(80, 14).Value = -2287.76923
$ws.Cells.Item(83, 8).Value = 282.44446
$ws.Cells.Item(83, 9).Value = 258.2
$ws.Cells.Item(83, 10).Value = 291.76923
$ws.Cells.Item(83, 11).Value = 1291
$ws.Cells.Item(83, 12).Value = 1458.84615
$ws.Cells.Item(83, 13).Value = 3701
$ws.Cells.Item(83, 14).Value = -11442.84615
$ws.Cells.Item(86, 8).Value = 2002.9231
$ws.Cells.Item(86, 9).Value = 1656.5483
$ws.Cells.Item(86, 10).Value = 2514.238
$ws.Cells.Item(86, 11).Value = 1656.5483
$ws.Cells.Item(86, 12).Value = 2514.238
$ws.Cells.Item(86, 13).Value = -533.5482999999999
$ws.Cells.Item(86, 14).Value = -4760.237999999999
$ws.Cells.Item(89, 8).Value = 2002.9231
$ws.Cells.Item(89, 9).Value = 1656.5483
$ws.Cells.Item(89, 10).Value = 2514.238
$ws.Cells.Item(89, 11).Value = 8282.7415
$ws.Cells.Item(89, 12).Value = 12571.19
$ws.Cells.Item(89, 13).Value = -2666.7415
$ws.Cells.Item(89, 14).Value = -23803.19
$ws.Cells.Item(134, 8).Value = 3352485.5
$ws.Cells.Item(134, 9).Value = 4279486.5
$ws.Cells.Item(134, 10).Value = 1019.46155
$ws.Cells.Item(134, 11).Value = 12838459.5
$ws.Cells.Item(134, 12).Value = 3058.38465
$ws.Cells.Item(134, 13).Value = -12835924.5
$ws.Cells.Item(134, 14).Value = -8128.38465
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2600689.8
$ws.Cells.Item(31, 9).Value = 952.71155
$ws.Cells.Item(31, 10).Value = 9715759
$ws.Cells.Item(31, 11).Value = 952.71155
$ws.Cells.Item(31, 12).Value = 9715759
$ws.Cells.Item(31, 13).Value = -657.71155
$ws.Cells.Item(31, 14).Value = -9716349
$ws.Cells.Item(34, 8).Value = 2600689.8
$ws.Cells.Item(34, 9).Value = 952.71155
$ws.Cells.Item(34, 10).Value = 9715759
$ws.Cells.Item(34, 11).Value = 952.71155
$ws.Cells.Item(34, 12).Value = 9715759
$ws.Cells.Item(34, 13).Value = -750.71155
$ws.Cells.Item(34, 14).Value = -9716163
$ws.Cells.Item(122, 8).Value = 1057.8235
$ws.Cells.Item(122, 9).Value = 1086.8182
$ws.Cells.Item(122, 10).Value = 1004.6667
$ws.Cells.Item(122, 11).Value = 3260.4546
$ws.Cells.Item(122, 12).Value = 3014.0001
$ws.Cells.Item(122, 13).Value = -810.4546
$ws.Cells.Item(122, 14).Value = -7914.0001
$ws.Cells.Item(132, 8).Value = 1481.0878
$ws.Cells.Item(132, 9).Value = 1461.0416
$ws.Cells.Item(132, 10).Value = 1588
$ws.Cells.Item(132, 11).Value = 4383.1248
$ws.Cells.Item(132, 12).Value = 4764
$ws.Cells.Item(132, 13).Value = -1853.1248
$ws.Cells.Item(132, 14).Value = -9824
$ws.Cells.Item(134, 8).Value = 1235.439
$ws.Cells.Item(134, 9).Value = 1271.6285
$ws.Cells.Item(134, 11).Value = 3814.8855
$ws.Cells.Item(134, 13).Value = -1279.8855
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 7353599
$ws.Cells.Item(5, 9).Value = 293.26666
$ws.Cells.Item(5, 10).Value = 13158840
$ws.Cells.Item(5, 11).Value = 879.79998
$ws.Cells.Item(5, 12).Value = 39476520
$ws.Cells.Item(5, 13).Value = -767.79998
$ws.Cells.Item(5, 14).Value = -39476744
$ws.Cells.Item(68, 8).Value = 1918.2812
$ws.Cells.Item(68, 9).Value = 901.5
$ws.Cells.Item(68, 10).Value = 2257.2083
$ws.Cells.Item(68, 11).Value = 2704.5
$ws.Cells.Item(68, 12).Value = 6771.624899999999
$ws.Cells.Item(68, 13).Value = -1893.5
$ws.Cells.Item(68, 14).Value = -8393.624899999999
$ws.Cells.Item(71, 8).Value = 1918.2812
$ws.Cells.Item(71, 9).Value = 901.5
$ws.Cells.Item(71, 10).Value = 2257.2083
$ws.Cells.Item(71, 11).Value = 8113.5
$ws.Cells.Item(71, 12).Value = 20314.8747
$ws.Cells.Item(71, 13).Value = -4057.5
$ws.Cells.Item(71, 14).Value = -28426.8747
$ws.Cells.Item(105, 8).Value = 1335831.6
$ws.Cells.Item(105, 10).Value = 1335831.6
$ws.Cells.Item(105, 12).Value = 4007494.8
$ws.Cells.Item(105, 14).Value = -4012736.8
$ws.Cells.Item(107, 8).Value = 50440.85
$ws.Cells.Item(107, 9).Value = 125320.75
$ws.Cells.Item(107, 10).Value = 31720.875
$ws.Cells.Item(107, 11).Value = 375962.25
$ws.Cells.Item(107, 12).Value = 95162.625
$ws.Cells.Item(107, 13).Value = -374042.25
$ws.Cells.Item(107, 14).Value = -99002.625
$ws.Cells.Item(113, 8).Value = 474.25882
$ws.Cells.Item(113, 9).Value = 468.27585
$ws.Cells.Item(113, 10).Value = 487.1111
$ws.Cells.Item(113, 11).Value = 1404.82755
$ws.Cells.Item(113, 12).Value = 1461.3333
$ws.Cells.Item(113, 13).Value = 765.17245
$ws.Cells.Item(113, 14).Value = -5801.3333
$ws.Cells.Item(122, 8).Value = 11161977
$ws.Cells.Item(122, 9).Value = 23810268
$ws.Cells.Item(122, 10).Value = 1324417.8
$ws.Cells.Item(122, 11).Value = 214292412
$ws.Cells.Item(122, 12).Value = 11919760.2
$ws.Cells.Item(122, 13).Value = -214289962
$ws.Cells.Item(122, 14).Value = -11924660.2
$ws.Cells.Item(132, 8).Value = 76924730
$ws.Cells.Item(132, 9).Value = 166667500
$ws.Cells.Item(132, 10).Value = 2352.8572
$ws.Cells.Item(132, 11).Value = 1500007500
$ws.Cells.Item(132, 12).Value = 21175.7148
$ws.Cells.Item(132, 13).Value = -1500004970
$ws.Cells.Item(132, 14).Value = -26235.7148
$ws.Cells.Item(135, 8).Value = 7353599
$ws.Cells.Item(135, 9).Value = 293.26666
$ws.Cells.Item(135, 10).Value = 13158840
$ws.Cells.Item(135, 11).Value = 2639.39994
$ws.Cells.Item(135, 12).Value = 118429560
$ws.Cells.Item(135, 13).Value = -104.3999400000002
$ws.Cells.Item(135, 14).Value = -118434630
$ws.Cells.Item(140, 8).Value = 1451.5555
$ws.Cells.Item(140, 9).Value = 709.4761999999999
$ws.Cells.Item(140, 11).Value = 2128.4286
$ws.Cells.Item(140, 13).Value = 3051.5714
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2047.4445
$ws.Cells.Item(113, 9).Value = 1519.8182
$ws.Cells.Item(113, 10).Value = 2876.5715
$ws.Cells.Item(113, 11).Value = 1519.8182
$ws.Cells.Item(113, 12).Value = 2876.5715
$ws.Cells.Item(113, 13).Value = 650.1818000000001
$ws.Cells.Item(113, 14).Value = -7216.5715
$ws.Cells.Item(132, 8).Value = 27029010
$ws.Cells.Item(132, 9).Value = 37039044
$ws.Cells.Item(132, 10).Value = 1911.4
$ws.Cells.Item(132, 11).Value = 111117132
$ws.Cells.Item(132, 12).Value = 5734.200000000001
$ws.Cells.Item(132, 13).Value = -111114602
$ws.Cells.Item(132, 14).Value = -10794.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2281.8333
$ws.Cells.Item(7, 9).Value = 2152
$ws.Cells.Item(7, 10).Value = 2346.75
$ws.Cells.Item(7, 11).Value = 2152
$ws.Cells.Item(7, 12).Value = 2346.75
$ws.Cells.Item(7, 13).Value = -2040
$ws.Cells.Item(7, 14).Value = -2570.75
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(40, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 2045
$ws.Cells.Item(122, 9).Value = 2045
$ws.Cells.Item(122, 11).Value = 6135
$ws.Cells.Item(122, 13).Value = -3685
$ws.Cells.Item(126, 8).Value = 2281.8333
$ws.Cells.Item(126, 9).Value = 2152
$ws.Cells.Item(126, 10).Value = 2346.75
$ws.Cells.Item(126, 11).Value = 6456
$ws.Cells.Item(126, 12).Value = 7040.25
$ws.Cells.Item(126, 13).Value = -3986
$ws.Cells.Item(126, 14).Value = -11980.25
$ws.Cells.Item(132, 8).Value = 3801.9023
$ws.Cells.Item(132, 9).Value = 3801.9023
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 11405.7069
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -8875.706900000001
$ws.Cells.Item(132, 14).ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 950.1667
$ws.Cells.Item(107, 9).Value = 302
$ws.Cells.Item(107, 10).Value = 1079.8
$ws.Cells.Item(107, 11).Value = 906
$ws.Cells.Item(107, 12).Value = 3239.4
$ws.Cells.Item(107, 14).Value = -7079.4
$ws.Cells.Item(107, 13).Value = 1014
$ws.Cells.Item(122, 8).Value = 1341.9131
$ws.Cells.Item(122, 9).Value = 1392.7
$ws.Cells.Item(122, 10).Value = 1003.3333
$ws.Cells.Item(122, 11).Value = 4178.1
$ws.Cells.Item(122, 12).Value = 3009.9999
$ws.Cells.Item(122, 13).Value = -1728.1
$ws.Cells.Item(122, 14).Value = -7909.9999
$ws.Cells.Item(126, 8).Value = 2749.9167
$ws.Cells.Item(126, 9).Value = 2749.9167
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 8249.750100000001
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -5779.750100000001
$ws.Cells.Item(126, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 7265552.5
$ws.Cells.Item(132, 9).Value = 7544980
$ws.Cells.Item(132, 10).Value = 425
$ws.Cells.Item(132, 11).Value = 22634940
$ws.Cells.Item(132, 12).Value = 1275
$ws.Cells.Item(132, 13).Value = -22632410
$ws.Cells.Item(132, 14).Value = -6335
$ws.Cells.Item(136, 8).Value = 2510554.8
$ws.Cells.Item(136, 9).Value = 5718.421
$ws.Cells.Item(136, 10).Value = 7520227
$ws.Cells.Item(136, 11).Value = 17155.263
$ws.Cells.Item(136, 12).Value = 22560681
$ws.Cells.Item(136, 13).Value = -14605.263
$ws.Cells.Item(136, 14).Value = -22565781
